$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ------------------------------------------------------------------
# 1) Merge the split "E" / "nd user should be able to download..."
#    runs into a single run with the full sentence text.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "End user should be able to download the generated article (UI should have download button/symbol).",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "End user should be able to download the generated article (UI should have download button/symbol).",
    2) | Out-Null

# ------------------------------------------------------------------
# Helper: append a brand-new paragraph after the current last one and
# replace its content with an explicit run of OOXML (keeps separate
# <w:r> elements distinct instead of letting identically-formatted
# runs coalesce).
# ------------------------------------------------------------------
function Add-RawParagraph([string]$innerXml) {
    $last = $d.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
    $p = $d.Paragraphs.Last
    $p.Range.Text = "PLACEHOLDER"
    $full = '<w:p ' + $wNs + '>' + $innerXml + '</w:p>'
    $p.Range.InsertXML($full)
}

$numPr3 = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>'
$dash = [char]0x2013

# ------------------------------------------------------------------
# 2) Append the new "Version-2" section at the end of the document.
# ------------------------------------------------------------------

# "Version-2:" (carries a lastRenderedPageBreak marker before the text)
Add-RawParagraph('<w:r><w:lastRenderedPageBreak/><w:t>Version-2:</w:t></w:r>')

# "Language - Telugu implementation"
Add-RawParagraph('<w:r><w:t>Language ' + $dash + ' Telugu implementation</w:t></w:r>')

# "Authorization Implementation:"
Add-RawParagraph('<w:r><w:t>Authorization Implementation:</w:t></w:r>')

# "Protected Route: /generate-article" + " - implemented in V0.1"
Add-RawParagraph($numPr3 + '<w:r><w:t>Protected Route: /generate-article</w:t></w:r><w:r><w:t xml:space="preserve"> ' + $dash + ' implemented in V0.1</w:t></w:r>')

# "User-Specific Article History: /article-history" + " " + "- implemented in V0.1"
Add-RawParagraph($numPr3 + '<w:r><w:t>User-Specific Article History: /article-history</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>' + $dash + ' implemented in V0.1</w:t></w:r>')

# "Download Article: /download/{article_id}" + " " + "- implemented in V0.1"
Add-RawParagraph($numPr3 + '<w:r><w:t>Download Article: /download/{</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>article_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>}</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>' + $dash + ' implemented in V0.1</w:t></w:r>')

# "Role-Based Access: /admin/users"
Add-RawParagraph($numPr3 + '<w:r><w:t>Role-Based Access: /admin/users</w:t></w:r>')

# "Token Expiry Handling"
Add-RawParagraph($numPr3 + '<w:r><w:t>Token Expiry Handling</w:t></w:r>')

# "Account Info Endpoint: /me"
Add-RawParagraph($numPr3 + '<w:r><w:t>Account Info Endpoint: /me</w:t></w:r>')
